$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date-of-birth values in column I (cells are formatted as dates already)
$ws.Range("I2").Value = (Get-Date -Year 1996 -Month 1 -Day 15 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("I3").Value = (Get-Date -Year 1981 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("I4").Value = (Get-Date -Year 1977 -Month 1 -Day 11 -Hour 0 -Minute 0 -Second 0).Date

# Update the active window view: scroll back to A1 (remove topLeftCell="D1")
# and move/ update the selection to I11
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("I11").Select()
